$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 44
$ws.Range("C6").Value = 44
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = "Plastic"
